$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.106.84'
$ws.Range('E2').Value = '  -1.43%  '
$ws.Range('D3').Value = '2.942.11'
$ws.Range('E3').Value = '  -2.27%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '376.50'
$ws.Range('E5').Value = '  -1.37%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '102.53'
$ws.Range('E6').Value = '  -4.30%  '
$ws.Range('E7').Value = '  -2.36%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -3.12%  '
$ws.Range('E10').Value = '  -3.05%  '
$ws.Range('E11').Value = '  -1.18%  '
$ws.Range('E12').Value = '  -1.53%  '
$ws.Range('D13').Value = '3.401.91'
$ws.Range('E13').Value = '  -2.41%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '17.96'
$ws.Range('E14').Value = '  -4.51%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.35'
$ws.Range('E15').Value = '  -2.92%  '
$ws.Range('D16').Value = '2.933.07'
$ws.Range('E16').Value = '  -2.40%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.974'
$ws.Range('E17').Value = '  -0.14%  '
$ws.Range('D18').Value = '51.047.33'
$ws.Range('E18').Value = '  -1.60%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.16'
$ws.Range('E19').Value = '  -6.71%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.13'
$ws.Range('E20').Value = '  -4.55%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.53'
$ws.Range('E21').Value = '  -4.98%  '
$ws.Range('D22').Value = '0.0₃0949'
$ws.Range('E22').Value = '  -1.61%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '263.14'
$ws.Range('E23').Value = '  -0.71%  '
$ws.Range('E24').Value = '  -1.32%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.86'
$ws.Range('E25').Value = '  +1.83%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.20'
$ws.Range('E26').Value = '  +9.16%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.75'
$ws.Range('E27').Value = '  +6.91%  '
$ws.Range('E28').Value = '  -3.58%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('E30').Value = '  +4.37%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '25.66'
$ws.Range('E31').Value = '  -2.33%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '9.85'
$ws.Range('E32').Value = '  -1.53%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '34.12'
$ws.Range('E33').Value = '  -3.06%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '50.66'
$ws.Range('E34').Value = '  -1.78%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0455'
$ws.Range('E35').Value = '  -0.48%  '
$ws.Range('E36').Value = '  -3.15%  '
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('E39').Value = '  -2.68%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '16.48'
$ws.Range('E40').Value = '  -6.78%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.115'
$ws.Range('E41').Value = '  -2.03%  '
$ws.Range('E42').Value = '  -5.03%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '121.43'
$ws.Range('E43').Value = '  -2.56%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '21.28'
$ws.Range('E44').Value = '  -5.43%  '
$ws.Range('E45').Value = '  -1.65%  '
$ws.Range('E46').Value = '  -2.74%  '
$ws.Range('E47').Value = '  -1.14%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.23'
$ws.Range('E48').Value = '  -2.33%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '2.007.94'
$ws.Range('E49').Value = '  -2.64%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0349'
$ws.Range('E50').Value = '  -2.39%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '5.02'
$ws.Range('E51').Value = '  -4.09%  '
